$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.140.66'
$ws.Range("E2").Value = '  +5.62%  '

$ws.Range("D3").Value = '1.920.97'
$ws.Range("E3").Value = '  +2.32%  '

$ws.Range("E4").Value = '  -1.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.14'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.24%  '

$ws.Range("E6").Value = '  -0.98%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5168'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.63%  '

$ws.Range("E8").Value = '  +2.94%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08458'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.77'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.124'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.76%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.32'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +9.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.341'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.89%  '

$ws.Range("D14").Value = '1.923.27'
$ws.Range("E14").Value = '  +2.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.358'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.11'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001116'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06719'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.23'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.055'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.01%  '

$ws.Range("D23").Value = '30.145.61'
$ws.Range("E23").Value = '  +5.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.26'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.200'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.63%  '

$ws.Range("D26").Value = '2.143.71'
$ws.Range("E26").Value = '  +2.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.67'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.16'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.458'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.95'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.078'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.64%  '

$ws.Range("E32").Value = '  +1.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.082'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.655'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.95%  '

$ws.Range("E35").Value = '  +2.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06591'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.69%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2214'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.44%  '

$ws.Range("E38").Value = '  +3.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6550'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.247'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.50'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6148'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.25'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.759'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.062'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.245'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.31%  '

$ws.Range("E49").Value = '  +2.68%  '

$ws.Range("E50").Value = '  +2.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.41'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.69%  '

# Row 39 <-> Row 40 data swap (index/A column stays the same)
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.211'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.60%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.010'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.49%  '

